$d = $word.ActiveDocument

# The "Requisitos" section ends with a paragraph reading
# "LOQ4057: Operações Unitárias III (Requisito fraco)". It used to be
# followed by a blank paragraph, then "Ver no Jupiter Salvar em pdf
# Salvar em docx", then the site footer "© 2020 . Contact:
# luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme
# under Creative Commons Attribution", then another blank paragraph and
# finally the page-break paragraph.
#
# The footer / "Ver no Jupiter" lines (plus the now-redundant blank
# paragraph between them and "LOQ4057: ...") are removed, leaving
# "LOQ4057: ..." directly followed by the remaining blank paragraph and
# the page break.

$paras = $d.Paragraphs
$count = $paras.Count

$startIndex = -1
$endIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $paras.Item($i).Range.Text
    if ($t -like "LOQ4057:*Requisito fraco)*") {
        $startIndex = $i
    }
    if ($t -like "*Contact: luizeleno@usp.br*Creative Commons Attribution*") {
        $endIndex = $i
    }
}

if ($startIndex -eq -1 -or $endIndex -eq -1 -or $endIndex -le $startIndex) {
    throw "Could not locate the paragraph range to delete (start=$startIndex end=$endIndex)"
}

$startPara = $paras.Item($startIndex)
$endPara = $paras.Item($endIndex)

# Delete from the end of the "LOQ4057: ..." paragraph (i.e. its
# paragraph mark) through the end of the "© 2020 ..." paragraph (i.e.
# its paragraph mark too), removing the three in-between paragraphs.
$rng = $d.Range($startPara.Range.End, $endPara.Range.End)
$rng.Delete()
